$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new header cells G1 (Elapsed Time) and H1 (CPU), matching F1's style ---
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# --- Fill in the new Elapsed Time / CPU columns for each data row ---
$ws.Range("G2").Value = 0.4788041146331428
$ws.Range("H2").Value = 0.997

$ws.Range("G3").Value = 0.4788041146331428
$ws.Range("H3").Value = 0.997

$ws.Range("G4").Value = 0.4788041146331428
$ws.Range("H4").Value = 0.997

$ws.Range("G5").Value = 0.4788041146331428
$ws.Range("H5").Value = 0.997

# --- Tiny floating point precision refresh on existing MSE/MAE values ---
$ws.Range("B3").Value = 0.1647362719327806
$ws.Range("D3").Value = 0.3215706312231009

$ws.Range("B4").Value = 0.2192987281846224
$ws.Range("D4").Value = 0.3777406617731509

$ws.Range("B5").Value = 0.3638702225807679
$ws.Range("D5").Value = 0.468209296615762
